$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, so numeric-looking strings
# like '6.280' or '1.572.31' keep their exact literal representation instead
# of being auto-coerced to a Number (which would drop trailing zeros etc.)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '22.475.19'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.572.31'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '291.39'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '0.3752'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = '49.89'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').Value = '0.07585'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '6.018'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('D15').Value = '6.964'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '1.568.76'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = '0.00001125'
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').Value = '91.21'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = '0.06741'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').Value = '6.280'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').Value = '16.46'
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').Value = '12.25'
$ws.Range('E23').Value = '  +2.05%  '
$ws.Range('D24').Value = '22.462.78'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '2.330'
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('D26').Value = '2.610'
$ws.Range('E26').Value = '  -5.94%  '
$ws.Range('D27').Value = '20.16'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').Value = '148.57'
$ws.Range('E28').Value = '  +2.44%  '
$ws.Range('D29').Value = '4.989'
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('D30').Value = '126.16'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '1.746.47'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '1.050'
$ws.Range('E32').Value = '  +3.31%  '
$ws.Range('D33').Value = '6.175'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('D35').Value = '9.893'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').Value = '0.08453'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('D37').Value = '1.384'
$ws.Range('E37').Value = '  +3.99%  '
$ws.Range('E38').Value = '  -3.25%  '
$ws.Range('D39').Value = '0.2298'
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('D40').Value = '0.06562'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Value = '5.489'
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').Value = '0.6306'
$ws.Range('E43').Value = '  -2.49%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '14.06'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = '0.5899'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('D48').Value = '2.102'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = '130.15'
$ws.Range('E49').Value = '  +3.67%  '
$ws.Range('D50').Value = '1.231'
$ws.Range('E50').Value = '  -5.35%  '
$ws.Range('D51').Value = '0.07335'
$ws.Range('E51').Value = '  -0.05%  '

# Restore default (unstyled) appearance for column D so no stray style index
# is left referenced on cells that originally had none.
$ws.Range('D2:D51').Style = 'Normal'
